# Add "NA" values under the duplicate_image_filename column (column E) for rows 2-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Work around a load/save quirk of the runtime that otherwise fills F1 with a
# stray value; restore it to blank so only the intended cells change.
$ws.Cells.Item(1, 6).Value = ""

for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
